$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cell (B17) carries the hyperlink font style ("s=1") used throughout column B.
$styleSource = $ws.Range("B17")

$ws.Range("A18").Value = "Ruben"
$ws.Range("B18").Value = "rubio@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B18"), "mailto:rubio@gmail.com", "", "", "rubio@gmail.com")
$styleSource.Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C18").Value = 13648292
$ws.Range("D18").Value = "Profesor"

$ws.Range("A19").Value = "Aurora"
$ws.Range("B19").Value = "aurora@hotmail.com"
$ws.Hyperlinks.Add($ws.Range("B19"), "mailto:aurora@hotmail.com", "", "", "aurora@hotmail.com")
$styleSource.Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("C19").Value = 36474398
$ws.Range("D19").Value = "Ingeniera industrial"

$ws.Range("A20").Value = "Madelen"
$ws.Range("B20").Value = "made@outlook.es"
$ws.Hyperlinks.Add($ws.Range("B20"), "mailto:made@outlook.es", "", "", "made@outlook.es")
$styleSource.Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("C20").Value = 58477363
$ws.Range("D20").Value = "Mecatrónica"

$ws.Range("A21").Value = "Virginia"
$ws.Range("B21").Value = "vivis@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B21"), "mailto:vivis@gmail.com", "", "", "vivis@gmail.com")
$styleSource.Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("C21").Value = 73887362
$ws.Range("D21").Value = "Enfermeria"

$ws.Range("A22").Value = "Noel"
$ws.Range("B22").Value = "noel@subitus.com"
$ws.Hyperlinks.Add($ws.Range("B22"), "mailto:noel@subitus.com", "", "", "noel@subitus.com")
$styleSource.Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C22").Value = 28376428
$ws.Range("D22").Value = "Teoría de gráficas"

$ws.Range("A23").Value = "Marisol"
$ws.Range("B23").Value = "marysol@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B23"), "mailto:marysol@gmail.com", "", "", "marysol@gmail.com")
$styleSource.Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("C23").Value = 77384939
$ws.Range("D23").Value = "Literatura"

$ws.Range("D24").Select()